$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its original "General" text-like display by forcing
# Text number format on cells whose new value would otherwise be parsed as a
# numeric literal by Excel (e.g. "230.08", "0.616"), matching the source data
# which stores these as plain text (inline strings) rather than numbers.

$ws.Range('D2').Value = '37.954.02'
$ws.Range('E2').Value = '  +2.16%  '

$ws.Range('D3').Value = '2.053.41'
$ws.Range('E3').Value = '  +1.40%  '

$ws.Range('E4').Value = '  +0.25%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '230.08'
$ws.Range('E5').Value = '  +1.23%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.616'
$ws.Range('E6').Value = '  +1.21%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.18'
$ws.Range('E7').Value = '  +5.48%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.387'
$ws.Range('E9').Value = '  +2.34%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0808'
$ws.Range('E10').Value = '  +2.65%  '

$ws.Range('E11').Value = '  +0.82%  '

$ws.Range('D12').Value = '2.355.54'
$ws.Range('E12').Value = '  +1.31%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.57'
$ws.Range('E13').Value = '  +2.45%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.72'
$ws.Range('E14').Value = '  +1.87%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.29'
$ws.Range('E15').Value = '  +2.51%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.751'
$ws.Range('E16').Value = '  +1.16%  '

$ws.Range('D17').Value = '2.051.58'
$ws.Range('E17').Value = '  +2.16%  '

$ws.Range('D18').Value = '37.890.98'
$ws.Range('E18').Value = '  +2.18%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.14'
$ws.Range('E19').Value = '  -2.49%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.66'
$ws.Range('E20').Value = '  +0.77%  '

$ws.Range('D21').Value = '0.0₃0831'
$ws.Range('E21').Value = '  +1.36%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.77'
$ws.Range('E22').Value = '  +0.23%  '

$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('E24').Value = '  +0.46%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.25'
$ws.Range('E25').Value = '  +2.27%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '166.66'
$ws.Range('E26').Value = '  +0.59%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  -0.04%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.133'
$ws.Range('E28').Value = '  +5.01%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.02'
$ws.Range('E29').Value = '  +1.18%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.36'
$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('E31').Value = '  +1.54%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  +0.27%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.58'
$ws.Range('E33').Value = '  +2.80%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0612'
$ws.Range('E34').Value = '  -0.76%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.05'
$ws.Range('E35').Value = '  +10.97%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.33'
$ws.Range('E36').Value = '  -0.90%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.98'
$ws.Range('E37').Value = '  +11.06%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.32'
$ws.Range('E38').Value = '  +4.96%  '

$ws.Range('E39').Value = '  -0.22%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0219'
$ws.Range('E40').Value = '  +0.85%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '97.67'
$ws.Range('E41').Value = '  +1.66%  '

$ws.Range('D42').Value = '1.485.74'
$ws.Range('E42').Value = '  +0.65%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0933'
$ws.Range('E44').Value = '  +1.42%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '16.66'
$ws.Range('E45').Value = '  +1.89%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.28'
$ws.Range('E46').Value = '  +19.37%  '

$ws.Range('E47').Value = '  -0.16%  '

$ws.Range('E48').Value = '  +0.08%  '

$ws.Range('E49').Value = '  +1.56%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.01'
$ws.Range('E50').Value = '  -2.77%  '

$ws.Range('D51').Value = '2.242.62'
$ws.Range('E51').Value = '  +1.41%  '
